# Testing user upload and report generation flow
# Update the "year" value for the student row from text "9D" to the numeric value 9,
# and leave the active selection on that cell (B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The year column (B) for row 2 was a text value "9D"; change it to the plain number 9.
$ws.Range("B2").Value = 9

# Move / leave the selection on B2 (matches the saved selection in the workbook).
$ws.Range("B2").Select()
